$wb = $excel.ActiveWorkbook

# --- Sheet 1: Battery_Data ---
$ws1 = $wb.Worksheets.Item("Battery_Data")
$ws1.Range("B2").Value = 532808.278649
$ws1.Range("B3").Value = 287716.47047046
$ws1.Range("B4").Value = 5754.329409409201
$ws1.Range("B5").Value = 72986.98354080001

# --- Sheet 2: Yearly BRC ---
$ws2 = $wb.Worksheets.Item("Yearly BRC")
$ws2.Range("B2").Value = 8893.01839920877
$ws2.Range("B3").Value = 8911.786229106507
$ws2.Range("B4").Value = 11024.66356877324
$ws2.Range("B5").Value = 11024.87703796586
$ws2.Range("B6").Value = 11024.87703796586
$ws2.Range("B7").Value = 11024.87703796586
$ws2.Range("B8").Value = 11021.40227621262
$ws2.Range("B9").Value = 14751.75318293285
$ws2.Range("B10").Value = 14747.52888169423
$ws2.Range("B11").Value = 14747.52888169423
$ws2.Range("B12").Value = 14747.52888169423
$ws2.Range("B13").Value = 14747.52888169423
$ws2.Range("B14").Value = 14747.52888169423
$ws2.Range("B15").Value = 14742.7237413523
$ws2.Range("B16").Value = 14029.60553099276
$ws2.Range("B17").Value = 14033.3818769272
$ws2.Range("B18").Value = 14033.3818769272
$ws2.Range("B19").Value = 14033.3818769272
$ws2.Range("B20").Value = 14033.38187692721
$ws2.Range("B21").Value = 14031.2710065673
